# "soem small docu fixes"
#  1. Refresh the auto-date footer placeholder (Slide Master + every Slide
#     Layout) from 14-Mar-18 to 27-Mar-18 - this is the cached display text
#     of the "Update automatically" date-and-time field.
#  2. Fix a typo on slide 2 ("LCE examples"):
#     lceGraphicGapentry -> lceGraphicTextEntry (and the run is no longer
#     flagged as a misspelling once corrected).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the date placeholder text everywhere it appears: the Slide
#    Master and each of its Custom Layouts. ppPlaceholderDate = 16.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "14-Mar-18") {
                $tr.Text = "27-Mar-18"
            }
        }
    }
}

Update-DatePlaceholder($p.SlideMaster)

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder($layouts.Item($li))
}

# ---------------------------------------------------------------------
# 2) Correct the typo on slide 2's "Content Placeholder 2" shape.
# ---------------------------------------------------------------------
$needle = "lceGraphicGapentry"
$replacement = "lceGraphicTextEntry"

$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $startPos = $full.IndexOf($needle)
        if ($startPos -ge 0) {
            $sub = $tr.Characters($startPos + 1, $needle.Length)
            $sub.Text = $replacement
        }
    }
}
